# Fixed single digit yards problem
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "yards" column (D) values that were missing a trailing digit
$updates = @{
    3  = 37
    4  = 59
    5  = 20
    6  = 22
    9  = 26
    10 = 74
    12 = 12
    13 = 87
    14 = 19
    15 = 18
    17 = 21
    21 = 41
    23 = 42
    25 = 23
    26 = 79
    27 = 72
    28 = 18
    31 = 10
    33 = 88
    34 = 68
    35 = 68
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}

# Restore the view so column C is the left-most visible column
$excel.ActiveWindow.ScrollColumn = 3
